$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying TPM recomputation dropped one sending/target-cluster combination
# (old row 5, FAPs/MuSCs -> Rspo2/Lgr6 -> MuSCs) and changed the remaining rows
# target clusters + all of the downstream expression-specificity numbers.
$ws.Rows.Item(5).Delete()

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo2"
$ws.Range("C2").Value = "Lgr6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1681403333333333
$ws.Range("H2").Value = 0.504421
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01182833333333333
$ws.Range("N2").Value = 0.035485
$ws.Range("O2").Value = 0.03045376408867423
$ws.Range("P2").Value = 0.03045376408867423
$ws.Range("Q2").Value = 0.001988819909444445
$ws.Range("R2").Value = 0.017899379185
$ws.Range("S2").Value = 0.03045376408867423
$ws.Range("T2").Value = 0.03045376408867423
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo2"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1681403333333333
$ws.Range("H3").Value = 0.504421
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2280253333333333
$ws.Range("N3").Value = 0.684076
$ws.Range("O3").Value = 0.5870843771374921
$ws.Range("P3").Value = 0.5870843771374921
$ws.Range("Q3").Value = 0.03834025555511111
$ws.Range("R3").Value = 0.345062299996
$ws.Range("S3").Value = 0.5870843771374921
$ws.Range("T3").Value = 0.5870843771374921
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo2"
$ws.Range("C4").Value = "Lgr6"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1681403333333333
$ws.Range("H4").Value = 0.504421
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1485493333333333
$ws.Range("N4").Value = 0.445648
$ws.Range("O4").Value = 0.3824618587738337
$ws.Range("P4").Value = 0.3824618587738337
$ws.Range("Q4").Value = 0.02497713442311111
$ws.Range("R4").Value = 0.224794209808
$ws.Range("S4").Value = 0.3824618587738337
$ws.Range("T4").Value = 0.3824618587738337

Write-Output "Updated Rspo2-Lgr6 LR-pair sheet with new TPM-based values"
